$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number need to be
# forced to text (NumberFormat "@") so Excel keeps the literal digit string
# instead of auto-converting it into a numeric value.

# Row 2
$ws.Range('D2').Value = '47.240.35'
$ws.Range('E2').Value = '  +4.46%  '

# Row 3
$ws.Range('D3').Value = '2.484.09'
$ws.Range('E3').Value = '  +2.10%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.43'
$ws.Range('E5').Value = '  +1.39%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.95'
$ws.Range('E6').Value = '  +1.43%  '

# Row 7
$ws.Range('E7').Value = '  +1.25%  '

# Row 8
$ws.Range('E8').Value = '  -0.05%  '

# Row 9
$ws.Range('E9').Value = '  +1.88%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.36'
$ws.Range('E10').Value = '  +4.43%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0810'
$ws.Range('E11').Value = '  +0.89%  '

# Row 12
$ws.Range('E12').Value = '  +0.24%  '

# Row 13
$ws.Range('E13').Value = '  -0.77%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.17'
$ws.Range('E14').Value = '  +2.72%  '

# Row 15
$ws.Range('D15').Value = '2.873.49'
$ws.Range('E15').Value = '  +2.09%  '

# Row 16
$ws.Range('D16').Value = '2.527.76'
$ws.Range('E16').Value = '  +3.70%  '

# Row 17
$ws.Range('E17').Value = '  +1.24%  '

# Row 18
$ws.Range('D18').Value = '47.143.41'
$ws.Range('E18').Value = '  +4.53%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.74'
$ws.Range('E19').Value = '  +3.93%  '

# Row 20
$ws.Range('E20').Value = '  +2.46%  '

# Row 21
$ws.Range('E21').Value = '  +0.77%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '70.59'
$ws.Range('E22').Value = '  +2.28%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '250.48'
$ws.Range('E23').Value = '  +3.01%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.37'
$ws.Range('E24').Value = '  +3.94%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.56'
$ws.Range('E25').Value = '  +1.93%  '

# Row 26
$ws.Range('E26').Value = '  +2.84%  '

# Row 27
$ws.Range('E27').Value = '  -0.08%  '

# Row 28
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.26'
$ws.Range('E28').Value = '  +3.71%  '

# Row 29
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.10'
$ws.Range('E29').Value = '  +5.80%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.34'
$ws.Range('E30').Value = '  +6.90%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.133'
$ws.Range('E31').Value = '  +5.10%  '

# Row 32
$ws.Range('E32').Value = '  +0.29%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.92'
$ws.Range('E33').Value = '  -2.05%  '

# Row 34
$ws.Range('E34').Value = '  +2.67%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0781'
$ws.Range('E35').Value = '  +1.94%  '

# Row 36
$ws.Range('E36').Value = '  +0.18%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.63'
$ws.Range('E37').Value = '  +3.46%  '

# Row 38
$ws.Range('E38').Value = '  +1.97%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.97'
$ws.Range('E39').Value = '  +4.17%  '

# Row 40
$ws.Range('E40').Value = '  +1.45%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '121.42'
$ws.Range('E41').Value = '  -2.27%  '

# Row 42
$ws.Range('E42').Value = '  +0.92%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.62'
$ws.Range('E43').Value = '  +2.29%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0295'
$ws.Range('E44').Value = '  +1.52%  '

# Row 45
$ws.Range('D45').Value = '1.949.25'
$ws.Range('E45').Value = '  +0.62%  '

# Row 46
$ws.Range('E46').Value = '  +1.48%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.09'
$ws.Range('E47').Value = '  +0.04%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.18'
$ws.Range('E48').Value = '  -0.74%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.78'
$ws.Range('E49').Value = '  +0.84%  '

# Row 50
$ws.Range('E50').Value = '  +13.53%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.71'
$ws.Range('E51').Value = '  +3.25%  '
